# Revised classroom banner image
# - Remove the empty Title/Subtitle placeholder shapes from slide 1
# - Reposition/resize the two banner picture shapes

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Remove the unused title & subtitle placeholders
$s.Shapes.Item("Title 1").Delete()
$s.Shapes.Item("Subtitle 2").Delete()

# Picture 3: move to its new position (size unchanged)
$pic3 = $s.Shapes.Item("Picture 3")
$pic3.Left = 113.9155
$pic3.Top = 108.0222

# Picture 4: move and resize (rotation unchanged)
$pic4 = $s.Shapes.Item("Picture 4")
$pic4.Left = 691.2098
$pic4.Top = 99.8501
$pic4.Width = 181.5751
$pic4.Height = 128.8597
